$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing number formats from the data rows so the new rows inherit
# the same styles (Entry/SL/TrailingSL/Target columns, PnL/MaxUp/MaxDown
# columns, and the MaxUp% column) already present in the sheet.
$fmtNum = $ws.Range("C2").NumberFormat        # columns C,D,E,F,G (s=2)
$fmtPnl = $ws.Range("H2").NumberFormat        # columns H,I,K (s=3)
$fmtPct = $ws.Range("J2").NumberFormat        # column J (s=4)

# ---- Row 6 ----
$ws.Range("A6").Value = "2025-10-23 22:44:23"
$ws.Range("B6").Value = "NSE:BANKNIFTY25OCT58400CE"

$ws.Range("C6").Value = 258.15
$ws.Range("C6").NumberFormat = $fmtNum

# D6 (LTP) stays blank for this trade, but keeps the same number format as
# the rest of the column.
$ws.Range("D6").NumberFormat = $fmtNum

$ws.Range("E6").Value = 243.15
$ws.Range("E6").NumberFormat = $fmtNum

$ws.Range("F6").Value = 243.15
$ws.Range("F6").NumberFormat = $fmtNum

$ws.Range("G6").Value = 288.15
$ws.Range("G6").NumberFormat = $fmtNum

# H6 (PnL) stays blank, but keeps the PnL number format.
$ws.Range("H6").NumberFormat = $fmtPnl

$ws.Range("I6").Value = 339.5000000000016
$ws.Range("I6").NumberFormat = $fmtPnl

$ws.Range("J6").Value = 3.75750532636066
$ws.Range("J6").NumberFormat = $fmtPct

$ws.Range("K6").Value = 61.25
$ws.Range("K6").NumberFormat = $fmtPnl

# ---- Row 7 ----
$ws.Range("A7").Value = "2025-10-31 10:31:45"
$ws.Range("B7").Value = "NSE:BANKNIFTY25NOV57900PE"

$ws.Range("C7").Value = 539.5
$ws.Range("C7").NumberFormat = $fmtNum

$ws.Range("D7").Value = 569.55
$ws.Range("D7").NumberFormat = $fmtNum

$ws.Range("E7").Value = 524.5
$ws.Range("E7").NumberFormat = $fmtNum

$ws.Range("F7").Value = 554.55
$ws.Range("F7").NumberFormat = $fmtNum

$ws.Range("G7").Value = 569.5
$ws.Range("G7").NumberFormat = $fmtNum

$ws.Range("H7").Value = 1051.749999999998
$ws.Range("H7").NumberFormat = $fmtPnl

$ws.Range("I7").Value = 1051.749999999998
$ws.Range("I7").NumberFormat = $fmtPnl

$ws.Range("J7").Value = 5.569972196478212
$ws.Range("J7").NumberFormat = $fmtPct

$ws.Range("K7").Value = -257.2500000000008
$ws.Range("K7").NumberFormat = $fmtPnl
